# Use 6.9% DM not 5.1% for digestate
$wb = $excel.ActiveWorkbook

$slurry = $wb.Worksheets.Item("Slurry")

# Digestate ("Afgasset biomasse") dry-matter rows (man.dm column C)
$slurry.Range("C4").Value = 6.9
$slurry.Range("C7").Value = 6.9

# Make the Slurry sheet the active tab (moves selection away from Climate)
# with C8 selected, matching the saved view state.
$slurry.Activate()
$slurry.Range("C8").Select()
